$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with the latest scraped figures. D values are stored as literal text in
# the source feed (e.g. "25.782.76", "1.000", "1.670"), so each Price cell
# is switched to Text format before its value is written - otherwise Excel
# would silently reinterpret the numeric-looking string as a number and
# drop the exact formatting (trailing zeros, dotted thousands groups, etc).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.782.76"
$ws.Range("E2").Value = "  -3.92%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.816.37"
$ws.Range("E3").Value = "  -2.94%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "277.87"
$ws.Range("E5").Value = "  -7.82%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5118"
$ws.Range("E7").Value = "  -4.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3513"
$ws.Range("E8").Value = "  -6.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.11"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06663"
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.96"
$ws.Range("E11").Value = "  -7.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8327"
$ws.Range("E12").Value = "  -6.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07909"
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.814.00"
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.077"
$ws.Range("E15").Value = "  -3.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.72"
$ws.Range("E16").Value = "  -5.87%  "
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("E18").Value = "  -3.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008026"
$ws.Range("E19").Value = "  -5.97%  "
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "25.839.11"
$ws.Range("E21").Value = "  -3.95%  "
$ws.Range("E22").Value = "  -4.82%  "
$ws.Range("E23").Value = "  -6.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.081"
$ws.Range("E24").Value = "  -4.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.63"
$ws.Range("E25").Value = "  -3.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.185"
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.670"
$ws.Range("E28").Value = "  -5.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.44"
$ws.Range("E29").Value = "  -4.23%  "
$ws.Range("E30").Value = "  -8.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.236"
$ws.Range("E31").Value = "  -7.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08842"
$ws.Range("E32").Value = "  -3.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04865"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7303"
$ws.Range("E34").Value = "  -8.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("E35").Value = "  -3.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.876"
$ws.Range("E36").Value = "  -3.93%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9997"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.359"
$ws.Range("E39").Value = "  -8.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5228"
$ws.Range("E40").Value = "  -13.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01848"
$ws.Range("E41").Value = "  -5.29%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9558"
$ws.Range("E42").Value = "  -10.95%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.206"
$ws.Range("E43").Value = "  -5.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.18"
$ws.Range("E44").Value = "  -4.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.059"
$ws.Range("E45").Value = "  -8.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4594"
$ws.Range("E47").Value = "  -10.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1368"
$ws.Range("E48").Value = "  -8.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.75"
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.202"
$ws.Range("E50").Value = "  -7.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.503"
$ws.Range("E51").Value = "  -7.70%  "
